$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old demo keyword list with the new "OPTIONAL SORTING AND FREE_SHIPMENT" list.
$ws.Range("A1").Value = "Aranacak Kelime"
$ws.Range("A2").Value = "Sleepy Bez 4 Numara"
$ws.Range("A3").Value = "Sleepy Islak Mendil"
$ws.Range("A4").Value = "Bebek Yürüteci"
$ws.Range("A5").Value = "Kitap"
$ws.Range("A6").Value = "Terlik"
$ws.Range("A7").Value = "Bilgisayar"
$ws.Range("A8").Value = "Laptop"
$ws.Range("A9").Value = "Televizyon LCD"
$ws.Range("A10").Value = "iphone X"
$ws.Range("A11").Value = "Samsung Galaxy S10"
$ws.Range("A12").Value = "Bardak"
$ws.Range("A13").Value = "Kazak"
$ws.Range("A14").Value = "Akülü Araba"

# Rows 8-14 pick up the same "free shipment" red-font highlight that already
# decorated the blank rows below them.
$ws.Range("A8:A14").Font.Color = $ws.Range("A16").Font.Color

# The old trailing placeholder row is gone now that row 14 is filled in.
$ws.Range("A15").Clear()

# Reflect the new working selection (the optional / highlighted block).
$ws.Range("A7:A14").Select()
